$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1285113333333333
$ws.Range("H2").Value = 0.385534
$ws.Range("I2").Value = 0.03749201237720504
$ws.Range("J2").Value = 0.03749201237720504
$ws.Range("M2").Value = 42.09975866666667
$ws.Range("N2").Value = 126.299276
$ws.Range("O2").Value = 0.3315552933456474
$ws.Range("P2").Value = 0.3315552933456474
$ws.Range("Q2").Value = 5.410296119264889
$ws.Range("R2").Value = 48.69266507338399
$ws.Range("S2").Value = 0.01243067516184286
$ws.Range("T2").Value = 0.01243067516184286
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1285113333333333
$ws.Range("H3").Value = 0.385534
$ws.Range("I3").Value = 0.03749201237720504
$ws.Range("J3").Value = 0.03749201237720504
$ws.Range("O3").Value = 0.4502223747274475
$ws.Range("P3").Value = 0.4502223747274475
$ws.Range("Q3").Value = 7.346697264925779
$ws.Range("R3").Value = 66.120275384332
$ws.Range("S3").Value = 0.01687974284577611
$ws.Range("T3").Value = 0.01687974284577611
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1285113333333333
$ws.Range("H4").Value = 0.385534
$ws.Range("I4").Value = 0.03749201237720504
$ws.Range("J4").Value = 0.03749201237720504
$ws.Range("M4").Value = 27.596267
$ws.Range("N4").Value = 82.78880100000001
$ws.Range("O4").Value = 0.2173335118824389
$ws.Range("P4").Value = 0.2173335118824389
$ws.Range("Q4").Value = 3.546433067192667
$ws.Range("R4").Value = 31.917897604734
$ws.Range("S4").Value = 0.008148270717477835
$ws.Range("T4").Value = 0.008148270717477838
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.1285113333333333
$ws.Range("H5").Value = 0.385534
$ws.Range("I5").Value = 0.03749201237720504
$ws.Range("J5").Value = 0.03749201237720504
$ws.Range("M5").Value = 0.1128593333333333
$ws.Range("N5").Value = 0.338578
$ws.Range("O5").Value = 0.0008888200444663087
$ws.Range("P5").Value = 0.0008888200444663087
$ws.Range("Q5").Value = 0.01450370340577778
$ws.Range("R5").Value = 0.130533330652
$ws.Range("S5").Value = 0.00003332365210823878
$ws.Range("T5").Value = 0.00003332365210823879
$ws.Range("I6").Value = 0.7552862722193517
$ws.Range("J6").Value = 0.755286272219352
$ws.Range("M6").Value = 42.09975866666667
$ws.Range("N6").Value = 126.299276
$ws.Range("O6").Value = 0.3315552933456474
$ws.Range("P6").Value = 0.3315552933456474
$ws.Range("Q6").Value = 108.9918126135813
$ws.Range("R6").Value = 980.9263135222319
$ws.Range("S6").Value = 0.2504191615456277
$ws.Range("T6").Value = 0.2504191615456277
$ws.Range("I7").Value = 0.7552862722193517
$ws.Range("J7").Value = 0.755286272219352
$ws.Range("O7").Value = 0.4502223747274475
$ws.Range("P7").Value = 0.4502223747274475
$ws.Range("S7").Value = 0.3400467790776379
$ws.Range("T7").Value = 0.340046779077638
$ws.Range("I8").Value = 0.7552862722193517
$ws.Range("J8").Value = 0.755286272219352
$ws.Range("M8").Value = 27.596267
$ws.Range("N8").Value = 82.78880100000001
$ws.Range("O8").Value = 0.2173335118824389
$ws.Range("P8").Value = 0.2173335118824389
$ws.Range("Q8").Value = 71.44381005869799
$ws.Range("R8").Value = 642.994290528282
$ws.Range("S8").Value = 0.1641490180180274
$ws.Range("T8").Value = 0.1641490180180275
$ws.Range("I9").Value = 0.7552862722193517
$ws.Range("J9").Value = 0.755286272219352
$ws.Range("M9").Value = 0.1128593333333333
$ws.Range("N9").Value = 0.338578
$ws.Range("O9").Value = 0.0008888200444663087
$ws.Range("P9").Value = 0.0008888200444663087
$ws.Range("Q9").Value = 0.2921808509106666
$ws.Range("R9").Value = 2.629627658196
$ws.Range("S9").Value = 0.0006713135780587968
$ws.Range("T9").Value = 0.0006713135780587969
$ws.Range("G10").Value = 0.692415
$ws.Range("H10").Value = 2.077245
$ws.Range("I10").Value = 0.2020057770533527
$ws.Range("J10").Value = 0.2020057770533527
$ws.Range("M10").Value = 42.09975866666667
$ws.Range("N10").Value = 126.299276
$ws.Range("O10").Value = 0.3315552933456474
$ws.Range("P10").Value = 0.3315552933456474
$ws.Range("Q10").Value = 29.15050439718
$ws.Range("R10").Value = 262.35453957462
$ws.Range("S10").Value = 0.06697608466843979
$ws.Range("T10").Value = 0.06697608466843981
$ws.Range("G11").Value = 0.692415
$ws.Range("H11").Value = 2.077245
$ws.Range("I11").Value = 0.2020057770533527
$ws.Range("J11").Value = 0.2020057770533527
$ws.Range("O11").Value = 0.4502223747274475
$ws.Range("P11").Value = 0.4502223747274475
$ws.Range("Q11").Value = 39.58377253389001
$ws.Range("R11").Value = 356.25395280501
$ws.Range("S11").Value = 0.09094752065362377
$ws.Range("T11").Value = 0.0909475206536238
$ws.Range("G12").Value = 0.692415
$ws.Range("H12").Value = 2.077245
$ws.Range("I12").Value = 0.2020057770533527
$ws.Range("J12").Value = 0.2020057770533527
$ws.Range("M12").Value = 27.596267
$ws.Range("N12").Value = 82.78880100000001
$ws.Range("O12").Value = 0.2173335118824389
$ws.Range("P12").Value = 0.2173335118824389
$ws.Range("Q12").Value = 19.108069214805
$ws.Range("R12").Value = 171.972622933245
$ws.Range("S12").Value = 0.04390262494754613
$ws.Range("T12").Value = 0.04390262494754613
$ws.Range("G13").Value = 0.692415
$ws.Range("H13").Value = 2.077245
$ws.Range("I13").Value = 0.2020057770533527
$ws.Range("J13").Value = 0.2020057770533527
$ws.Range("M13").Value = 0.1128593333333333
$ws.Range("N13").Value = 0.338578
$ws.Range("O13").Value = 0.0008888200444663087
$ws.Range("P13").Value = 0.0008888200444663087
$ws.Range("Q13").Value = 0.07814549528999999
$ws.Range("R13").Value = 0.70330945761
$ws.Range("S13").Value = 0.0001795467837430122
$ws.Range("T13").Value = 0.0001795467837430122
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.01787866666666667
$ws.Range("H14").Value = 0.053636
$ws.Range("I14").Value = 0.005215938350090445
$ws.Range("J14").Value = 0.005215938350090446
$ws.Range("M14").Value = 42.09975866666667
$ws.Range("N14").Value = 126.299276
$ws.Range("O14").Value = 0.3315552933456474
$ws.Range("P14").Value = 0.3315552933456474
$ws.Range("Q14").Value = 0.7526875519484445
$ws.Range("R14").Value = 6.774187967536
$ws.Range("S14").Value = 0.00172937196973705
$ws.Range("T14").Value = 0.00172937196973705
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.01787866666666667
$ws.Range("H15").Value = 0.053636
$ws.Range("I15").Value = 0.005215938350090445
$ws.Range("J15").Value = 0.005215938350090446
$ws.Range("O15").Value = 0.4502223747274475
$ws.Range("P15").Value = 0.4502223747274475
$ws.Range("Q15").Value = 1.022082240480889
$ws.Range("R15").Value = 9.198740164328001
$ws.Range("S15").Value = 0.002348332150409685
$ws.Range("T15").Value = 0.002348332150409685
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.01787866666666667
$ws.Range("H16").Value = 0.053636
$ws.Range("I16").Value = 0.005215938350090445
$ws.Range("J16").Value = 0.005215938350090446
$ws.Range("M16").Value = 27.596267
$ws.Range("N16").Value = 82.78880100000001
$ws.Range("O16").Value = 0.2173335118824389
$ws.Range("P16").Value = 0.2173335118824389
$ws.Range("Q16").Value = 0.4933844589373334
$ws.Range("R16").Value = 4.440460130436001
$ws.Range("S16").Value = 0.00113359819938745
$ws.Range("T16").Value = 0.001133598199387451
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.01787866666666667
$ws.Range("H17").Value = 0.053636
$ws.Range("I17").Value = 0.005215938350090445
$ws.Range("J17").Value = 0.005215938350090446
$ws.Range("M17").Value = 0.1128593333333333
$ws.Range("N17").Value = 0.338578
$ws.Range("O17").Value = 0.0008888200444663087
$ws.Range("P17").Value = 0.0008888200444663087
$ws.Range("Q17").Value = 0.002017774400888889
$ws.Range("R17").Value = 0.018159969608
$ws.Range("S17").Value = 0.000004636030556260915
$ws.Range("T17").Value = 0.000004636030556260915
